$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "xd"
$ws.Range("A1").Font.ThemeColor = 1
